$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 301 (the post about "古いデバイスを交換して新品をゲットしよう"),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(301).Delete()
